# Insert a new data row at row 618 (pushing the existing rows 618-659 down
# to 619-660) and populate it with the new entry:
#   2026/01/12, 月, 4, 201
#
# Using Rows.Item(618).Insert() shifts everything below down by one row,
# exactly like Excel's native "Insert Sheet Rows" command, and keeps the
# sheet's used-range dimension (A1:D659 -> A1:D660) in sync automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(618).Insert()

# Column A holds a date-look-alike string (e.g. "2026/01/12") that must stay
# literal text (matching every other row in the column), not get silently
# reinterpreted as a date serial by the smart-entry parser. Prefixing with an
# apostrophe forces literal text; resetting the Style back to "Normal"
# afterwards drops the quote-prefix formatting Excel would otherwise tack on,
# so the new cell ends up styled identically to its neighbours.
$ws.Range("A618").Value = "'2026/01/12"
$ws.Range("A618").Style = "Normal"

$ws.Range("B618").Value = "月"
$ws.Range("C618").Value = 4
$ws.Range("D618").Value = 201
